# Update ILR data / dashboard text (dataText sheet)
#
# - Rows 14-17, column B: refresh the "data period" note text from the
#   AY22/23 wording to the AY23/24 wording.
# - Rows 14-17, column D: refresh the ILR source-link text to point at the
#   new Explore Education Statistics data-catalogue dataset URL.
# - Move the active selection/viewport from E8 to C18 to reflect where the
#   author was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

$newPeriodNote = "AY23/24 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on"
$newIlrLink = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"

foreach ($r in 14..17) {
    $ws.Cells.Item($r, 2).Value = $newPeriodNote   # column B
    $ws.Cells.Item($r, 4).Value = $newIlrLink       # column D
}

# Reflect the new scroll position / selection (previously topLeftCell A8,
# selection E8) as topLeftCell A18, selection C18.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
